# Insert two new weekly data rows (Primera / Segunda) for Apio right after
# the existing row 394, shifting all subsequent rows down by two. This
# grows the used range from A1:R514 to A1:R516.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 395 (old rows 395..514 shift to 397..516).
$ws.Rows.Item(395).Resize(2).Insert()

# --- New row 395 ("Primera") ---
$ws.Range("A395").Value = 8
$ws.Range("B395").Value = "Terminal La Palmera de La Serena"
$ws.Range("C395").Value = "Coquimbo"
$ws.Range("D395").Value = 44876
$ws.Range("E395").Value = 4
$ws.Range("F395").Value = 100112017
$ws.Range("G395").Value = "Apio"
$ws.Range("H395").Value = "Americana (o)"
$ws.Range("I395").Value = "Primera"
$ws.Range("J395").Value = 1800
$ws.Range("K395").Value = 9000
$ws.Range("L395").Value = 10000
$ws.Range("M395").Value = 9500
$ws.Range("N395").Value = "$/docena de matas"
$ws.Range("O395").Value = "Provincia del Elquí"
$ws.Range("P395").Value = 1583
$ws.Range("Q395").Value = 6
$ws.Range("R395").Value = "Hortaliza"

# --- New row 396 ("Segunda") ---
$ws.Range("A396").Value = 8
$ws.Range("B396").Value = "Terminal La Palmera de La Serena"
$ws.Range("C396").Value = "Coquimbo"
$ws.Range("D396").Value = 44876
$ws.Range("E396").Value = 4
$ws.Range("F396").Value = 100112017
$ws.Range("G396").Value = "Apio"
$ws.Range("H396").Value = "Americana (o)"
$ws.Range("I396").Value = "Segunda"
$ws.Range("J396").Value = 1300
$ws.Range("K396").Value = 7000
$ws.Range("L396").Value = 8000
$ws.Range("M396").Value = 7500
$ws.Range("N396").Value = "$/docena de matas"
$ws.Range("O396").Value = "Provincia del Elquí"
$ws.Range("P396").Value = 1250
$ws.Range("Q396").Value = 6
$ws.Range("R396").Value = "Hortaliza"
